# fix parameterization, run strategy and influence experiments
#
# Rows 23 (CV SALTS management zones), 26 (CDFA) and 28 (Regional water
# management groups) already use the "red / negative influence" style in
# columns C:F, they just need the -1 values filled in (style untouched).
#
# Row 25 (water quality coalitions) needs both: the -1 values AND the
# same "red / negative influence" formatting that C23:F23 already use,
# so we copy that formatting across before writing the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: pick up the style used by the other "-1" cells (e.g. C23:F23)
# and apply it to C25:F25, then fill in the values.
$styleSource = $ws.Range("C23")
$styleSource.Copy()
$ws.Range("C25:F25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C25:F25").Value = -1

# --- Rows 23, 26, 28: formatting already correct, just set the values.
$ws.Range("C23:F23").Value = -1
$ws.Range("C26:F26").Value = -1
$ws.Range("C28:F28").Value = -1

# --- Update the active selection / view state to match the saved session.
$ws.Range("D21").Select()
